$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "Step 2" row text (C3/D3) to the clarified wording
$ws.Range("C3").Value = "Step 2: Login as a user with the appropiate role"
$ws.Range("D3").Value = "I am redirected to the dashboard"

# Fill in previously empty rows 4 and 5 with the new, more granular steps
$ws.Range("C4").Value = " Step 3: From the dashboard go to the ""Create Employee"""
$ws.Range("D4").Value = "I am redirected to the Create Employee page"

$ws.Range("C5").Value = "Step 4:  Fill out the proper employee information and submit"
$ws.Range("D5").Value = "A new employee is added to the database"

# Match the wrapped/top-aligned formatting used by the other data rows
$ws.Range("C4:D5").WrapText = $true
$ws.Range("C4:D5").VerticalAlignment = -4160

# Update selection to reflect the new active cell
$ws.Range("C5").Select()
$excel.ActiveWindow.ScrollRow = 2
